# "changes to new branch"
# Appends three new paragraphs ("AWS files", "New branch") each preceded
# by a blank paragraph, plus a trailing blank paragraph, right after the
# existing "README" paragraph (i.e. just before the section properties).

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-EmptyParagraphAtEnd {
    # Grow the document by one paragraph at the very end, then hand back
    # a Range scoped to just that new (still empty) paragraph.
    $last = $d.Paragraphs.Last.Range
    $last.InsertParagraphAfter()
    return $d.Paragraphs.Last.Range
}

function Set-ParagraphPlainText($paraRange, [string]$text) {
    # Range.InsertXML replaces the contents of the range it's called on,
    # so targeting the freshly-added (empty) paragraph's own Range lets
    # us stamp its final contents without disturbing its neighbours.
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xml = "<w:p $wNs><w:r><w:t>$escaped</w:t></w:r></w:p>"
    $paraRange.InsertXML($xml)
}

function Set-ParagraphEmpty($paraRange) {
    $xml = "<w:p $wNs/>"
    $paraRange.InsertXML($xml)
}

Set-ParagraphEmpty (New-EmptyParagraphAtEnd)
Set-ParagraphPlainText (New-EmptyParagraphAtEnd) "AWS files"
Set-ParagraphEmpty (New-EmptyParagraphAtEnd)
Set-ParagraphPlainText (New-EmptyParagraphAtEnd) "New branch"
Set-ParagraphEmpty (New-EmptyParagraphAtEnd)
